$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 8.467854666666666
$ws.Range("H2").Value2 = 25.403564
$ws.Range("I2").Value2 = 0.1523462492674013
$ws.Range("J2").Value2 = 0.1523462492674013
$ws.Range("M2").Value2 = 14.129345
$ws.Range("N2").Value2 = 42.388035
$ws.Range("O2").Value2 = 0.3414817166893976
$ws.Range("P2").Value2 = 0.3414817166893976
$ws.Range("Q2").Value2 = 119.6452399951933
$ws.Range("R2").Value2 = 1076.80715995674
$ws.Range("S2").Value2 = 0.05202345873102307
$ws.Range("T2").Value2 = 0.05202345873102308

# Row 3
$ws.Range("G3").Value2 = 8.467854666666666
$ws.Range("H3").Value2 = 25.403564
$ws.Range("I3").Value2 = 0.1523462492674013
$ws.Range("J3").Value2 = 0.1523462492674013
$ws.Range("O3").Value2 = 0.3553528814026711
$ws.Range("P3").Value2 = 0.3553528814026711
$ws.Range("Q3").Value2 = 124.5052917930529
$ws.Range("R3").Value2 = 1120.547626137476
$ws.Range("S3").Value2 = 0.05413667864806062
$ws.Range("T3").Value2 = 0.05413667864806063

# Row 4
$ws.Range("G4").Value2 = 8.467854666666666
$ws.Range("H4").Value2 = 25.403564
$ws.Range("I4").Value2 = 0.1523462492674013
$ws.Range("J4").Value2 = 0.1523462492674013
$ws.Range("O4").Value2 = 0.3031654019079313
$ws.Range("P4").Value2 = 0.3031654019079312
$ws.Range("Q4").Value2 = 106.2203201423693
$ws.Range("R4").Value2 = 955.982881281324
$ws.Range("S4").Value2 = 0.0461861118883176
$ws.Range("T4").Value2 = 0.0461861118883176

# Row 5
$ws.Range("I5").Value2 = 0.1669927598427297
$ws.Range("J5").Value2 = 0.1669927598427297
$ws.Range("M5").Value2 = 14.129345
$ws.Range("N5").Value2 = 42.388035
$ws.Range("O5").Value2 = 0.3414817166893976
$ws.Range("P5").Value2 = 0.3414817166893976
$ws.Range("Q5").Value2 = 131.147887952095
$ws.Range("R5").Value2 = 1180.330991568855
$ws.Range("S5").Value2 = 0.05702497430579564
$ws.Range("T5").Value2 = 0.05702497430579564

# Row 6
$ws.Range("I6").Value2 = 0.1669927598427297
$ws.Range("J6").Value2 = 0.1669927598427297
$ws.Range("O6").Value2 = 0.3553528814026711
$ws.Range("P6").Value2 = 0.3553528814026711
$ws.Range("S6").Value2 = 0.05934135838349829
$ws.Range("T6").Value2 = 0.05934135838349829

# Row 7
$ws.Range("I7").Value2 = 0.1669927598427297
$ws.Range("J7").Value2 = 0.1669927598427297
$ws.Range("O7").Value2 = 0.3031654019079313
$ws.Range("P7").Value2 = 0.3031654019079312
$ws.Range("S7").Value2 = 0.05062642715343581
$ws.Range("T7").Value2 = 0.0506264271534358

# Row 8
$ws.Range("I8").Value2 = 0.680660990889869
$ws.Range("J8").Value2 = 0.680660990889869
$ws.Range("M8").Value2 = 14.129345
$ws.Range("N8").Value2 = 42.388035
$ws.Range("O8").Value2 = 0.3414817166893976
$ws.Range("P8").Value2 = 0.3414817166893976
$ws.Range("Q8").Value2 = 534.557614657405
$ws.Range("R8").Value2 = 4811.018531916646
$ws.Range("S8").Value2 = 0.2324332836525789
$ws.Range("T8").Value2 = 0.2324332836525789

# Row 9
$ws.Range("I9").Value2 = 0.680660990889869
$ws.Range("J9").Value2 = 0.680660990889869
$ws.Range("O9").Value2 = 0.3553528814026711
$ws.Range("P9").Value2 = 0.3553528814026711
$ws.Range("S9").Value2 = 0.2418748443711122
$ws.Range("T9").Value2 = 0.2418748443711122

# Row 10
$ws.Range("I10").Value2 = 0.680660990889869
$ws.Range("J10").Value2 = 0.680660990889869
$ws.Range("O10").Value2 = 0.3031654019079313
$ws.Range("P10").Value2 = 0.3031654019079312
$ws.Range("S10").Value2 = 0.2063528628661779
$ws.Range("T10").Value2 = 0.2063528628661778

